# Updated symbol list on Thu Dec 15 04:58:48 UTC 2022 with GitHub Actions
#
# Refreshes the "Price" column (D) on Sheet1 with newly scraped quotes.
# The column stores values as plain text (not numbers), so each new value
# is written with a leading apostrophe to force Excel to keep it as text
# (preserving things like trailing zeros, e.g. "3.900" or "0.0001500").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceUpdates = [ordered]@{
    2  = "265.03"
    3  = "22.58"
    4  = "6.281"
    5  = "0.06139"
    6  = "3.593"
    7  = "6.665"
    8  = "1.342"
    9  = "0.8277"
    11 = "0.1585"
    12 = "0.08230"
    13 = "0.03413"
    14 = "0.03125"
    15 = "0.09246"
    16 = "3.900"
    17 = "0.001716"
    18 = "0.04879"
    19 = "0.006240"
    20 = "0.005269"
    22 = "0.0001500"
    23 = "3.764"
    24 = "2.315"
    26 = "0.1237"
    27 = "0.0002681"
    40 = "0.04601"
    41 = "0.006958"
    42 = "0.1136"
    43 = "0.003401"
    44 = "0.01081"
    45 = "0.00006171"
    46 = "0.00000000750"
    47 = "0.7784"
    48 = "0.1933"
    49 = "0.00002101"
    50 = "0.01240"
}

foreach ($row in $priceUpdates.Keys) {
    $ws.Range("D$row").Value = "'" + $priceUpdates[$row]
}
